# Remove the existing "_GoBack" bookmark (it currently wraps the Builder
# image paragraph). Deleting it removes both the <w:bookmarkStart> and
# <w:bookmarkEnd> markers.
$d = $word.ActiveDocument
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the paragraph that contains the lone "Builder" heading text.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Builder", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse to the end of the found text ("Builder") so new text is appended
# immediately after it, inside the same paragraph.
$searchRange.Collapse(0)
$insertStart = $searchRange.Start

$newText = ", permet de créer l’instance Restaurant"
# Append a temporary placeholder character after the new text so that, while
# we are adding the bookmark, its position is never the very last character
# slot of the paragraph (immediately before the paragraph mark) -- the COM
# host mis-handles Bookmarks.Add exactly at that boundary. We delete the
# placeholder right after the bookmark is created.
$placeholder = "X"
$searchRange.InsertAfter($newText + $placeholder)

# The appended text should not be bold, unlike the "Builder" run; restore
# the rest of the run formatting (Calibri / bCs / iCs) by only toggling Bold
# off over the inserted text range.
$newTextRange = $d.Range($insertStart, $insertStart + $newText.Length)
$newTextRange.Bold = 0

# Re-create the "_GoBack" bookmark, collapsed, right after the new text
# (i.e. immediately before the paragraph mark).
$bookmarkPos = $insertStart + $newText.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the temporary placeholder character now that the bookmark is safely
# anchored in place.
$placeholderRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$placeholderRange.Delete()
